$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates extracted from the target diff: (cell, new value)
$updates = @(
    [PSCustomObject]@{ Cell = 'D2'; Value = '24.451.13' }
    [PSCustomObject]@{ Cell = 'E2'; Value = '  +8.84%  ' }
    [PSCustomObject]@{ Cell = 'D3'; Value = '1.681.56' }
    [PSCustomObject]@{ Cell = 'E3'; Value = '  +5.05%  ' }
    [PSCustomObject]@{ Cell = 'D4'; Value = '1.002' }
    [PSCustomObject]@{ Cell = 'E4'; Value = '  +0.07%  ' }
    [PSCustomObject]@{ Cell = 'D5'; Value = '306.85' }
    [PSCustomObject]@{ Cell = 'E5'; Value = '  +4.93%  ' }
    [PSCustomObject]@{ Cell = 'D6'; Value = '0.9974' }
    [PSCustomObject]@{ Cell = 'E6'; Value = '  +0.46%  ' }
    [PSCustomObject]@{ Cell = 'D7'; Value = '0.3719' }
    [PSCustomObject]@{ Cell = 'E7'; Value = '  +0.70%  ' }
    [PSCustomObject]@{ Cell = 'D8'; Value = '0.3458' }
    [PSCustomObject]@{ Cell = 'E8'; Value = '  +1.90%  ' }
    [PSCustomObject]@{ Cell = 'D9'; Value = '48.00' }
    [PSCustomObject]@{ Cell = 'E9'; Value = '  +12.39%  ' }
    [PSCustomObject]@{ Cell = 'D10'; Value = '1.185' }
    [PSCustomObject]@{ Cell = 'E10'; Value = '  +3.88%  ' }
    [PSCustomObject]@{ Cell = 'D11'; Value = '0.07279' }
    [PSCustomObject]@{ Cell = 'E11'; Value = '  +3.35%  ' }
    [PSCustomObject]@{ Cell = 'D12'; Value = '0.9995' }
    [PSCustomObject]@{ Cell = 'E12'; Value = '  +0.31%  ' }
    [PSCustomObject]@{ Cell = 'D13'; Value = '20.47' }
    [PSCustomObject]@{ Cell = 'E13'; Value = '  +3.69%  ' }
    [PSCustomObject]@{ Cell = 'D14'; Value = '6.145' }
    [PSCustomObject]@{ Cell = 'E14'; Value = '  +3.63%  ' }
    [PSCustomObject]@{ Cell = 'D15'; Value = '6.759' }
    [PSCustomObject]@{ Cell = 'E15'; Value = '  +1.93%  ' }
    [PSCustomObject]@{ Cell = 'D16'; Value = '1.676.82' }
    [PSCustomObject]@{ Cell = 'E16'; Value = '  +4.73%  ' }
    [PSCustomObject]@{ Cell = 'E17'; Value = '  +2.50%  ' }
    [PSCustomObject]@{ Cell = 'D18'; Value = '0.9970' }
    [PSCustomObject]@{ Cell = 'E18'; Value = '  +0.46%  ' }
    [PSCustomObject]@{ Cell = 'D19'; Value = '0.06726' }
    [PSCustomObject]@{ Cell = 'E19'; Value = '  +1.70%  ' }
    [PSCustomObject]@{ Cell = 'D20'; Value = '81.52' }
    [PSCustomObject]@{ Cell = 'E20'; Value = '  +4.40%  ' }
    [PSCustomObject]@{ Cell = 'D21'; Value = '16.49' }
    [PSCustomObject]@{ Cell = 'E21'; Value = '  +2.30%  ' }
    [PSCustomObject]@{ Cell = 'E22'; Value = '  +1.38%  ' }
    [PSCustomObject]@{ Cell = 'E23'; Value = '  +1.80%  ' }
    [PSCustomObject]@{ Cell = 'D24'; Value = '24.385.76' }
    [PSCustomObject]@{ Cell = 'E24'; Value = '  +8.82%  ' }
    [PSCustomObject]@{ Cell = 'D25'; Value = '2.433' }
    [PSCustomObject]@{ Cell = 'E25'; Value = '  +1.08%  ' }
    [PSCustomObject]@{ Cell = 'D26'; Value = '2.677' }
    [PSCustomObject]@{ Cell = 'E26'; Value = '  +6.62%  ' }
    [PSCustomObject]@{ Cell = 'D27'; Value = '3.362' }
    [PSCustomObject]@{ Cell = 'E27'; Value = '  -11.63%  ' }
    [PSCustomObject]@{ Cell = 'D28'; Value = '153.06' }
    [PSCustomObject]@{ Cell = 'E28'; Value = '  +2.16%  ' }
    [PSCustomObject]@{ Cell = 'D29'; Value = '19.60' }
    [PSCustomObject]@{ Cell = 'E29'; Value = '  +0.27%  ' }
    [PSCustomObject]@{ Cell = 'D30'; Value = '1.861.14' }
    [PSCustomObject]@{ Cell = 'E30'; Value = '  +4.88%  ' }
    [PSCustomObject]@{ Cell = 'D31'; Value = '126.95' }
    [PSCustomObject]@{ Cell = 'E31'; Value = '  +5.60%  ' }
    [PSCustomObject]@{ Cell = 'D32'; Value = '6.340' }
    [PSCustomObject]@{ Cell = 'E32'; Value = '  +5.16%  ' }
    [PSCustomObject]@{ Cell = 'D33'; Value = '4.027' }
    [PSCustomObject]@{ Cell = 'E33'; Value = '  -4.38%  ' }
    [PSCustomObject]@{ Cell = 'D34'; Value = '0.9731' }
    [PSCustomObject]@{ Cell = 'E34'; Value = '  +2.03%  ' }
    [PSCustomObject]@{ Cell = 'D35'; Value = '1.720' }
    [PSCustomObject]@{ Cell = 'E35'; Value = '  +5.96%  ' }
    [PSCustomObject]@{ Cell = 'D36'; Value = '0.08481' }
    [PSCustomObject]@{ Cell = 'E36'; Value = '  +2.92%  ' }
    [PSCustomObject]@{ Cell = 'B37'; Value = 'Aptos' }
    [PSCustomObject]@{ Cell = 'C37'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    [PSCustomObject]@{ Cell = 'D37'; Value = '12.48' }
    [PSCustomObject]@{ Cell = 'E37'; Value = '  +5.39%  ' }
    [PSCustomObject]@{ Cell = 'B38'; Value = 'Hedera' }
    [PSCustomObject]@{ Cell = 'C38'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    [PSCustomObject]@{ Cell = 'D38'; Value = '0.06514' }
    [PSCustomObject]@{ Cell = 'E38'; Value = '  +6.58%  ' }
    [PSCustomObject]@{ Cell = 'B39'; Value = 'FraxShare' }
    [PSCustomObject]@{ Cell = 'C39'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    [PSCustomObject]@{ Cell = 'D39'; Value = '9.015' }
    [PSCustomObject]@{ Cell = 'E39'; Value = '  +3.94%  ' }
    [PSCustomObject]@{ Cell = 'D40'; Value = '5.366' }
    [PSCustomObject]@{ Cell = 'E40'; Value = '  +1.61%  ' }
    [PSCustomObject]@{ Cell = 'E41'; Value = '  +5.90%  ' }
    [PSCustomObject]@{ Cell = 'D42'; Value = '1.266' }
    [PSCustomObject]@{ Cell = 'E42'; Value = '  +1.44%  ' }
    [PSCustomObject]@{ Cell = 'D43'; Value = '0.2112' }
    [PSCustomObject]@{ Cell = 'E43'; Value = '  +4.27%  ' }
    [PSCustomObject]@{ Cell = 'D44'; Value = '0.6197' }
    [PSCustomObject]@{ Cell = 'E44'; Value = '  +4.89%  ' }
    [PSCustomObject]@{ Cell = 'D45'; Value = '0.9976' }
    [PSCustomObject]@{ Cell = 'E45'; Value = '  +0.42%  ' }
    [PSCustomObject]@{ Cell = 'D46'; Value = '3.782' }
    [PSCustomObject]@{ Cell = 'E46'; Value = '  +2.62%  ' }
    [PSCustomObject]@{ Cell = 'E47'; Value = '  +4.49%  ' }
    [PSCustomObject]@{ Cell = 'D48'; Value = '13.02' }
    [PSCustomObject]@{ Cell = 'E48'; Value = '  -0.67%  ' }
    [PSCustomObject]@{ Cell = 'D49'; Value = '127.12' }
    [PSCustomObject]@{ Cell = 'E49'; Value = '  +0.51%  ' }
    [PSCustomObject]@{ Cell = 'D50'; Value = '2.035' }
    [PSCustomObject]@{ Cell = 'E50'; Value = '  +3.60%  ' }
    [PSCustomObject]@{ Cell = 'D51'; Value = '0.07229' }
    [PSCustomObject]@{ Cell = 'E51'; Value = '  +6.02%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell -match '^D') {
        # Price column: values look like plain decimals (e.g. 48.00, 0.9970) which
        # Excel would otherwise silently reinterpret as numbers and strip trailing
        # zeros from. Force literal text, then restore the default (Normal) style
        # so no stray per-cell number format lingers on the cell.
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
